$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.6811
$ws.Range("C7").Value = -11.9744
$ws.Range("B8").Value = 4.689600000000001
$ws.Range("A12").Value = -22.76610000000002
$ws.Range("B12").Value = 6.1213
$ws.Range("B14").Value = 8.411000000000007
$ws.Range("C19").Value = -13.45179999999999
$ws.Range("E19").Value = 13.4306
$ws.Range("C21").Value = -13.0117
$ws.Range("B22").Value = 4.758800000000005
$ws.Range("C24").Value = -11.48699999999999
